# ---------------------------------------------------------------------------
# Edit script: rewrites several paragraphs of the Introduction / Education
# sections, fixes a couple of typos/spell-check markers in the table, and
# removes a duplicated "Related Work" heading near the end of the document.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function New-Pkg([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# standard run-properties block used throughout the body text
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr>'
$rPrSuper = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:vertAlign w:val="superscript"/><w:lang w:val="en-US"/></w:rPr>'

function Replace-ParagraphContent([string]$anchorText, [string]$newBodyRunsXml) {
    $rng = $d.Content
    $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $p = $rng.Paragraphs(1)
    $contentRng = $d.Range($p.Range.Start, $p.Range.End - 1)
    $contentRng.InsertXML((New-Pkg $newBodyRunsXml))
}

# ---------------------------------------------------------------------------
# 1. Introduction paragraph: split "Augmented Reality (AR) and Virtual
#    Reality (VR) are two exciting..." into 3 runs with reworded text.
# ---------------------------------------------------------------------------
$body1 = "<w:p><w:r>$rPr<w:t xml:space=`"preserve`">Augmented Reality (AR) and </w:t></w:r>" +
         "<w:r>$rPr<w:t>Virtual Reality (VR)</w:t></w:r>" +
         "<w:r>$rPr<w:t xml:space=`"preserve`"> are two great technologies that have evolved the way people interact with the digital world. AR is the technology to add digital elements to the real world around us, while VR is the technology to create immersive digital environments that may or may not replicate the real world in sense and form. Both technologies use special equipment such as headset or glasses to bring these experiences to life.</w:t></w:r></w:p>"
Replace-ParagraphContent "Augmented Reality (AR) and Virtual Reality (VR) are two exciting technologies" $body1

# ---------------------------------------------------------------------------
# 2. "Augmented reality (AR) is an enhanced version..." paragraph rewrite.
# ---------------------------------------------------------------------------
$body2 = "<w:p>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Augmented reality (AR) can be considered as an enhanced version of the real world, achieved through the utilization of digital information and projecting to show in semblance with the real-world environment. At the core of </w:t></w:r>" +
    "<w:proofErr w:type=`"gramStart`"/>" +
    "<w:r>$rPr<w:t>it</w:t></w:r>" +
    "<w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> AR is used to add features to real world using computer generated digital information, whether the features are for utilization or display they all add </w:t></w:r>" +
    "<w:r>$rPr<w:t>new information</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> to our interpretation of the real-world environment that elevates our experience towards the said environment. AR can be utilized through multiple senses including visual, auditory, </w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">touch etc. AR technology uses computer hardware and </w:t></w:r>" +
    "<w:r>$rPr<w:lastRenderedPageBreak/><w:t xml:space=`"preserve`">software for example apps, consoles, projections etc., to combine digital information with the </w:t></w:r>" +
    "<w:r>$rPr<w:t>real-world</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> environment. </w:t></w:r>" +
    "</w:p>"
Replace-ParagraphContent "Augmented reality (AR) is an enhanced version of the real world" $body2

# ---------------------------------------------------------------------------
# 3. "AR is a growing trend..." -> reworded single sentence.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("AR is a growing trend among companies developing metaverse solutions, particularly in mobile computing and business applications.", $true, $false, $false, $false, $false, $true, 1, $false, "It is growing trend among companies developing metaverse implementations such as mobile computing and business applications, to use AR as in their implementation.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. "Augmented reality either makes visual changes..." + "...gaming, product
#    visualization..." (2 runs incl. page break) merged into a single run.
# ---------------------------------------------------------------------------
$body4 = "<w:p><w:r>$rPr<w:t>Various fields which have applications for AR include gaming, product visualization, marketing campaigns, architecture and home design, education etc.</w:t></w:r></w:p>"
Replace-ParagraphContent "Augmented reality either makes visual changes to a natural environment" $body4

# ---------------------------------------------------------------------------
# 5. VR training paragraph ("An important area of application for VR
#    systems...") gets replaced with new "Education has always been..." text.
# ---------------------------------------------------------------------------
$body5 = "<w:p>" +
    "<w:r>$rPr<w:t>Education has always been an important area of application for VR systems so that students could efficiently train for real life activities. The appeal of simulations was exactly that they can provide training that could stand equally with practice with real systems, adding the reduced cost and greater safety along with that made it tree worth growing for companies. This was particularly found in the case of military training, where the first significant commercial application of simulations was seen</w:t></w:r>" +
    "<w:r>$rPr<w:t>, it was pilot training simulators used in the 2</w:t></w:r>" +
    "<w:r>$rPrSuper<w:t>nd</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> World War. Such simulators relied on visual-motion feedback to augment the sensation of flying while the user is seated in a closed mechanical system placed on the ground.</w:t></w:r>" +
    "</w:p>"
Replace-ParagraphContent "An important area of application for VR systems has always been training" $body5

# ---------------------------------------------------------------------------
# 6. "As highlighted above, AR/VR technology enhances..." paragraph content
#    is entirely removed, leaving an empty paragraph (pPr kept).
# ---------------------------------------------------------------------------
$rng6 = $d.Content
$rng6.Find.Execute("As highlighted above, AR/VR technology enhances", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p6 = $rng6.Paragraphs(1)
$content6 = $d.Range($p6.Range.Start, $p6.Range.End - 1)
$content6.Text = ""

# ---------------------------------------------------------------------------
# 7. Insert two brand-new paragraphs before "This Research paper will
#    further explore..." recreating (with edits) the old VR-training and
#    "As highlighted above" paragraphs.
# ---------------------------------------------------------------------------
$rng7 = $d.Content
$rng7.Find.Execute("This Research paper will further explore", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p7 = $rng7.Paragraphs(1)
$insertPos = $p7.Range.Start

$p7.Range.InsertParagraphBefore()
$p7.Range.InsertParagraphBefore()

$bodyA = "<w:p><w:r>$rPr<w:t xml:space=`"preserve`">An important area of application for VR systems has always been training for real-life activities. The appeal of simulations is that they can provide training equal or nearly equal to practice with real systems, but at reduced cost and with greater safety. This is particularly the case for military training, and the first significant application of commercial simulators was </w:t></w:r>" +
         "<w:r>$rPr<w:lastRenderedPageBreak/><w:t>pilot training during World War II. Flight simulators rely on visual and motion feedback to augment the sensation of flying while seated in a closed mechanical system on the ground.</w:t></w:r></w:p>"
$targetA = $d.Range($insertPos, $insertPos)
$targetA.InsertXML((New-Pkg $bodyA))

$bodyB = "<w:p>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">As highlighted above, AR/VR technology enhances the learning process of people in various domains and fields through simulations which may put people in various positions by simulating various situations or examples of situations, thereby enhancing experiential learning. This point can be further enunciated with examples </w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">of other fields as well, simulated surgical training and skill improvement is known to be one of the most revolutionary uses of AR/VR technology in the healthcare industry, </w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">in countries like </w:t></w:r>" +
    "<w:r>$rPr<w:t>Netherlands</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">, and </w:t></w:r>" +
    "<w:proofErr w:type=`"gramStart`"/>" +
    "<w:r>$rPr<w:t>United kingdom</w:t></w:r>" +
    "<w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> the </w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">police departments </w:t></w:r>" +
    "<w:r>$rPr<w:t>is being trained using AR/VR technology to train them for various emergency scenarios, companies like Walmart have also implemented VR training programs to train employees in various departments etc.</w:t></w:r>" +
    "</w:p>"

$rngB = $d.Content
$rngB.Find.Execute("This Research paper will further explore", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pB = $rngB.Paragraphs(1)
$targetB = $d.Range($pB.Previous().Range.Start, $pB.Previous().Range.Start)
$targetB.InsertXML((New-Pkg $bodyB))

# ---------------------------------------------------------------------------
# 8. "This Research paper will further explore..." paragraph: merge the
#    trailing two runs (dropping the mid-sentence page break) into one.
# ---------------------------------------------------------------------------
$body8 = "<w:p>" +
    "<w:r>$rPr<w:t>This Research paper will further explore and discuss such innovations,</w:t></w:r>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> their degree of benefit and effectiveness, as well as further and upcoming implementations. This study will mainly employ existing literature, case studies and empirical research to elucidate the role of AR/VR technology in shaping the future of skill acquisition and lifelong learning.</w:t></w:r>" +
    "</w:p>"
Replace-ParagraphContent "This Research paper will further explore and discuss such innovations" $body8

# ---------------------------------------------------------------------------
# 9. Table: "Patrice Labedan" -> split with spell-check proofErr markers.
# ---------------------------------------------------------------------------
$rng9 = $d.Content
$rng9.Find.Execute("Patrice Labedan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$xml9 = "<w:p>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">Patrice </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r>$rPr<w:t>Labedan</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "</w:p>"
$rng9.InsertXML((New-Pkg $xml9))

# ---------------------------------------------------------------------------
# 10. Table: "VR's adaptability and  potential for personalized learning
#     experiences" -> split with grammar-check proofErr markers.
# ---------------------------------------------------------------------------
$rng10 = $d.Content
$rng10.Find.Execute("VR's adaptability and  potential for personalized learning experiences", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$xml10 = "<w:p>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`">VR's adaptability </w:t></w:r>" +
    "<w:proofErr w:type=`"gramStart`"/>" +
    "<w:r>$rPr<w:t>and  potential</w:t></w:r>" +
    "<w:proofErr w:type=`"gramEnd`"/>" +
    "<w:r>$rPr<w:t xml:space=`"preserve`"> for personalized learning experiences</w:t></w:r>" +
    "</w:p>"
$rng10.InsertXML((New-Pkg $xml10))

# ---------------------------------------------------------------------------
# 11. Remove the duplicated "Related Work" heading paragraph near the end of
#     the document (keeps the page-break paragraph before it and the blank
#     paragraph after it).
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item($d.Tables.Count)
$rng11 = $d.Range($tbl.Range.End, $d.Content.End)
$rng11.Find.Execute("Related Work", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$p11 = $rng11.Paragraphs(1)
$delRng = $d.Range($p11.Range.Start, $p11.Range.End)
$delRng.Delete()

Write-Host "All edits applied."
